$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original Text format so
# numeric-looking values like "1.00" / "227.83" are not reinterpreted
# as numbers (which would drop formatting such as trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.797.21'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.104.73'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.83'
$ws.Range("E5").Value = '  -0.07%  '

$ws.Range("E6").Value = '  +0.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.22'
$ws.Range("E7").Value = '  +2.43%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  +2.32%  '

$ws.Range("E10").Value = '  +0.89%  '

$ws.Range("E11").Value = '  -0.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.82'
$ws.Range("E12").Value = '  +6.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.416.59'
$ws.Range("E13").Value = '  +0.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.02'
$ws.Range("E14").Value = '  -1.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.808'
$ws.Range("E15").Value = '  +2.80%  '

$ws.Range("E16").Value = '  +1.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.089.48'
$ws.Range("E17").Value = '  +0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.806.85'
$ws.Range("E18").Value = '  +1.30%  '

$ws.Range("E19").Value = '  +1.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.60'
$ws.Range("E20").Value = '  +0.25%  '

$ws.Range("E21").Value = '  +1.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.27'
$ws.Range("E22").Value = '  +1.25%  '

$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  -3.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  +0.63%  '

$ws.Range("E26").Value = '  +2.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.14'
$ws.Range("E27").Value = '  +1.40%  '

$ws.Range("E28").Value = '  +1.18%  '

$ws.Range("E29").Value = '  +3.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.34'
$ws.Range("E30").Value = '  +1.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.58'
$ws.Range("E31").Value = '  +10.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.57'
$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.19'
$ws.Range("E34").Value = '  +12.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.75'
$ws.Range("E35").Value = '  -1.05%  '

$ws.Range("E36").Value = '  +2.06%  '

$ws.Range("E37").Value = '  +0.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.52'
$ws.Range("E38").Value = '  -0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.12'
$ws.Range("E40").Value = '  -1.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.24'
$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("E42").Value = '  +3.52%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.525.45'
$ws.Range("E43").Value = '  -0.82%  '

$ws.Range("E44").Value = '  +8.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  +0.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0917'
$ws.Range("E46").Value = '  -1.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.74'
$ws.Range("E47").Value = '  +0.69%  '

$ws.Range("E48").Value = '  +4.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.15'
$ws.Range("E49").Value = '  +0.53%  '

$ws.Range("E50").Value = '  -0.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.302.77'
$ws.Range("E51").Value = '  +0.85%  '
